$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample number text "E7420" -> "E7420L" for all rows (column G, rows 2-27)
$ws.Range("G2:G27").Value = "E7420L"

# Replace the "accuracy check" formula cells (=FALSE()) in column H with plain boolean FALSE values
$ws.Range("H2:H27").Value = $False

# Update the active selection to the accuracy-check column (H2:H27)
$ws.Range("H2:H27").Select()
